$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.274.17"
$ws.Range("E2").Value = "  +3.86%  "
$ws.Range("D3").Value = "1.812.36"
$ws.Range("E3").Value = "  +4.79%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "'329.33"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "'0.4461"
$ws.Range("E7").Value = "  +5.84%  "
$ws.Range("D8").Value = "'0.3709"
$ws.Range("D9").Value = "'44.83"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'0.07715"
$ws.Range("E10").Value = "  +4.54%  "
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").Value = "'0.9996"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "'22.13"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").Value = "'6.305"
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("D15").Value = "'7.577"
$ws.Range("E15").Value = "  +6.59%  "
$ws.Range("D16").Value = "1.847.24"
$ws.Range("E16").Value = "  +6.75%  "
$ws.Range("D17").Value = "'92.98"
$ws.Range("E17").Value = "  +7.39%  "
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("D19").Value = "'0.06550"
$ws.Range("E19").Value = "  +10.11%  "
$ws.Range("D20").Value = "'0.9999"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'17.51"
$ws.Range("E21").Value = "  +4.63%  "
$ws.Range("D22").Value = "'6.232"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D23").Value = "28.321.95"
$ws.Range("E23").Value = "  +3.88%  "
$ws.Range("D24").Value = "'11.70"
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("D25").Value = "'2.154"
$ws.Range("E25").Value = "  -9.82%  "
$ws.Range("D26").Value = "'20.78"
$ws.Range("E26").Value = "  +3.98%  "
$ws.Range("D27").Value = "'156.04"
$ws.Range("D28").Value = "2.036.77"
$ws.Range("E28").Value = "  +5.76%  "
$ws.Range("D29").Value = "'2.315"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").Value = "'128.33"
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("D31").Value = "'1.200"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").Value = "'5.911"
$ws.Range("E32").Value = "  +6.09%  "
$ws.Range("D33").Value = "'0.09258"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "'3.657"
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("D35").Value = "'13.07"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").Value = "'0.02359"
$ws.Range("E36").Value = "  +5.62%  "
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("D38").Value = "'5.188"
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("D39").Value = "'0.06237"
$ws.Range("E39").Value = "  +3.37%  "
$ws.Range("D40").Value = "'0.6586"
$ws.Range("E40").Value = "  +4.05%  "
$ws.Range("D41").Value = "'1.200"
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("D42").Value = "'8.161"
$ws.Range("E42").Value = "  +3.74%  "
$ws.Range("D43").Value = "'0.9995"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "'1.404"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'13.93"
$ws.Range("E45").Value = "  +4.02%  "
$ws.Range("D46").Value = "'0.6096"
$ws.Range("E46").Value = "  +5.36%  "
$ws.Range("D47").Value = "'3.770"
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("D48").Value = "'127.09"
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("D49").Value = "'2.038"
$ws.Range("E49").Value = "  +5.15%  "
$ws.Range("D50").Value = "'1.158"
$ws.Range("D51").Value = "'0.06985"
$ws.Range("E51").Value = "  +2.64%  "
